# ---------------------------------------------------------------------------
# Edit script: splits the two "fe_frontend last_version" rows out of
# rights_and_functions (sheet1) into a brand-new worksheet
# "noch_manuel_generierung_umsetze" (sheet3), adjusting the rights_and_functions
# table accordingly and adding descriptive header comments to the new sheet.
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("rights_and_functions")

# ---------------------------------------------------------------------------
# 1) Create the new worksheet as the last sheet in the workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "noch_manuel_generierung_umsetze"

# ---------------------------------------------------------------------------
# 2) Build the header block of the new sheet (mirrors the
#    rights_and_functions header at rows 22/23).
# ---------------------------------------------------------------------------
$ws3.Range("A1").Value2 = $ws1.Range("A22").Value2
$ws3.Range("A1").Style  = $ws1.Range("A22").Style

"B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q" | ForEach-Object {
    $ws3.Range("$_`1").Style = $ws1.Range("A4").Style
}

"A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q" | ForEach-Object {
    $ws3.Range("$_`2").Value2 = $ws1.Range("$_`23").Value2
    $ws3.Range("$_`2").Style  = $ws1.Range("$_`23").Style
}

# ---------------------------------------------------------------------------
# 3) Copy the two rows being relocated (old sheet1 rows 64 & 67) into the
#    new sheet as rows 3 & 4, keeping their original values/styles.
# ---------------------------------------------------------------------------
"B","C","D","E","G","H","I","J","N","O" | ForEach-Object {
    $ws3.Range("$_`3").Value2 = $ws1.Range("$_`64").Value2
    $ws3.Range("$_`3").Style  = $ws1.Range("$_`64").Style

    $ws3.Range("$_`4").Value2 = $ws1.Range("$_`67").Value2
    $ws3.Range("$_`4").Style  = $ws1.Range("$_`67").Style
}

# Row 3 gains a leading path label in column A; row 4's A cell stays blank
# but still carries the "A1-style" padding, matching the rest of the row.
$ws3.Range("A3").Value2 = $ws1.Range("A6").Value2
$ws3.Range("A3").Style  = $ws1.Range("A4").Style
$ws3.Range("A4").Style  = $ws1.Range("A4").Style

# Padding cells (blank, but explicitly styled) so the row matches the table
# layout used elsewhere in the workbook.
"F","K","L","M" | ForEach-Object {
    $ws3.Range("$_`3").Style = $ws1.Range("A4").Style
    $ws3.Range("$_`4").Style = $ws1.Range("A4").Style
}

# The manual/"noch_manuel_generierung_umsetze" rows use a dedicated suffix
# (new shared string) instead of the old "_last_version" one.
$ws3.Range("H3").Value2 = "_fe_last_version"
$ws3.Range("H4").Value2 = "_fe_last_version"

$ws3.Range("M7").Select()

# ---------------------------------------------------------------------------
# 4) Add the descriptive header comments to the new sheet (D2 / H2),
#    matching the ones already present on rights_and_functions (D23 / H23).
# ---------------------------------------------------------------------------
$ws3.Range("D2").AddComment("Autor:`nGRANT TRIGGER`nGRANT USAGE ON SCHEMA`nGRANT USAGE ON seq") | Out-Null
$ws3.Range("H2").AddComment("Autor:`nRawdaten = varchar oder kein Eintrag dann Datentypen") | Out-Null

# ---------------------------------------------------------------------------
# 5) Remove the two relocated rows from rights_and_functions. Delete the
#    lower-numbered row last so indices of the still-to-delete row don't
#    shift out from under us.
# ---------------------------------------------------------------------------
$ws1.Rows.Item(67).Delete()
$ws1.Rows.Item(64).Delete()

# ---------------------------------------------------------------------------
# 6) The two comments that lived below the deleted rows need to move up by
#    two rows (K71 -> K69, K78 -> K76); re-create them at their new address.
# ---------------------------------------------------------------------------
$oldK71 = $ws1.Range("K71").Comment.Text()
$ws1.Range("K71").Comment.Delete()
$ws1.Range("K69").AddComment($oldK71) | Out-Null

$oldK78 = $ws1.Range("K78").Comment.Text()
$ws1.Range("K78").Comment.Delete()
$ws1.Range("K76").AddComment($oldK78) | Out-Null

# ---------------------------------------------------------------------------
# 7) Restore the view state: rights_and_functions keeps the active tab,
#    with B65 selected; the new sheet should have M7 selected (already set
#    above, before we switched away from it).
# ---------------------------------------------------------------------------
$ws1.Range("B65").Select()
$ws1.Select()
